$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Sheet1 -> Details)
$ws.Name = "Details"

# Header row
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"

# Data rows
$ws.Range("A2").Value = "Tom"
$ws.Range("B2").Value = "Tom123"
$ws.Range("A3").Value = "John"
$ws.Range("B3").Value = "John123"

# Build the combined fill+border look in an out-of-the-way scratch cell first,
# so the "fill+border" cellXf is the first custom style registered (matches
# the header's eventual style slot), then stamp it onto the header via
# copy/paste-special so no throw-away intermediate style is minted.
$scratch = $ws.Range("Z1")
$scratch.Interior.Color = 49407
$scratch.Borders.LineStyle = 1

$scratch.Copy()
$headerRange = $ws.Range("A1:B1")
$headerRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false
$scratch.Clear()

# Data rows just get the plain box border (second distinct custom style).
$dataRange = $ws.Range("A2:B3")
$dataRange.Borders.LineStyle = 1

# Match the recorded selection from the saved workbook.
$ws.Range("E35").Select()
